$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# Row 17
$ws_ALC.Range("H17").Value = 893.9524
$ws_ALC.Range("I17").Value = 757
$ws_ALC.Range("J17").Value = 996.6667
$ws_ALC.Range("K17").Value = 2271
$ws_ALC.Range("L17").Value = 2990.0001
$ws_ALC.Range("M17").Value = -2103
$ws_ALC.Range("N17").Value = -3326.0001

# Row 33
$ws_ALC.Range("H33").Value = 6075.9414
$ws_ALC.Range("I33").Value = 97.09999999999999
$ws_ALC.Range("J33").Value = 14617.143
$ws_ALC.Range("K33").Value = 97.09999999999999
$ws_ALC.Range("L33").Value = 14617.143
$ws_ALC.Range("M33").Value = 131.9
$ws_ALC.Range("N33").Value = -15075.143

# Row 41
$ws_ALC.Range("H41").Value = 796
$ws_ALC.Range("I41").Value = 908.5
$ws_ALC.Range("J41").Value = 458.5
$ws_ALC.Range("K41").Value = 908.5
$ws_ALC.Range("L41").Value = 458.5
$ws_ALC.Range("M41").Value = -468.5
$ws_ALC.Range("N41").Value = -1338.5

# Row 43
$ws_ALC.Range("H43").Value = 2582.8333
$ws_ALC.Range("I43").Value = 0
$ws_ALC.Range("J43").Value = 2582.8333
$ws_ALC.Range("K43").Value = 0
$ws_ALC.Range("L43").ClearContents()
$ws_ALC.Range("M43").Value = 2582.8333
$ws_ALC.Range("N43").Value = -2720.8333

# Row 112
$ws_ALC.Range("H112").Value = 3620.75
$ws_ALC.Range("I112").Value = 640
$ws_ALC.Range("J112").Value = 4405.1577
$ws_ALC.Range("K112").Value = 1920
$ws_ALC.Range("L112").Value = 13215.4731
$ws_ALC.Range("M112").Value = -812
$ws_ALC.Range("N112").Value = -15431.4731

# Row 132
$ws_ALC.Range("H132").Value = 1574.3422
$ws_ALC.Range("I132").Value = 1372.8485
$ws_ALC.Range("J132").Value = 2904.2
$ws_ALC.Range("K132").Value = 4118.5455
$ws_ALC.Range("L132").Value = 8712.599999999999
$ws_ALC.Range("M132").Value = -1588.5455
$ws_ALC.Range("N132").Value = -13772.6

# Row 137
$ws_ALC.Range("H137").Value = 1547.4828
$ws_ALC.Range("I137").Value = 1186
$ws_ALC.Range("J137").Value = 2139
$ws_ALC.Range("K137").Value = 3558
$ws_ALC.Range("L137").Value = 6417
$ws_ALC.Range("M137").Value = -1008
$ws_ALC.Range("N137").Value = -11517

# Row 138
$ws_ALC.Range("H138").Value = 4182.6284
$ws_ALC.Range("I138").Value = 4158.1
$ws_ALC.Range("J138").Value = 4186.717
$ws_ALC.Range("K138").Value = 12474.3
$ws_ALC.Range("L138").Value = 12560.151
$ws_ALC.Range("M138").Value = -7334.300000000001
$ws_ALC.Range("N138").Value = -22840.151

$ws_ARM = $wb.Worksheets.Item("ARM")
# Row 74
$ws_ARM.Range("H74").Value = 875.9655
$ws_ARM.Range("I74").Value = 866.2593000000001
$ws_ARM.Range("J74").Value = 1007
$ws_ARM.Range("K74").Value = 866.2593000000001
$ws_ARM.Range("L74").Value = 1007
$ws_ARM.Range("M74").Value = 7.740699999999947
$ws_ARM.Range("N74").Value = -2755

# Row 77
$ws_ARM.Range("H77").Value = 875.9655
$ws_ARM.Range("I77").Value = 866.2593000000001
$ws_ARM.Range("J77").Value = 1007
$ws_ARM.Range("K77").Value = 4331.2965
$ws_ARM.Range("L77").Value = 5035
$ws_ARM.Range("M77").Value = 36.70349999999962
$ws_ARM.Range("N77").Value = -13771

$ws_CRP = $wb.Worksheets.Item("CRP")
# Row 31
$ws_CRP.Range("H31").Value = 17675.086
$ws_CRP.Range("I31").Value = 0
$ws_CRP.Range("J31").Value = 17675.086
$ws_CRP.Range("K31").Value = 0
$ws_CRP.Range("L31").ClearContents()
$ws_CRP.Range("M31").Value = 17675.086
$ws_CRP.Range("N31").Value = -18265.086

# Row 34
$ws_CRP.Range("H34").Value = 17675.086
$ws_CRP.Range("I34").Value = 0
$ws_CRP.Range("J34").Value = 17675.086
$ws_CRP.Range("K34").Value = 0
$ws_CRP.Range("L34").ClearContents()
$ws_CRP.Range("M34").Value = 17675.086
$ws_CRP.Range("N34").Value = -18079.086

# Row 55
$ws_CRP.Range("H55").Value = 6000
$ws_CRP.Range("I55").Value = 6000
$ws_CRP.Range("J55").Value = 0
$ws_CRP.Range("K55").Value = 6000
$ws_CRP.Range("L55").Value = 0
$ws_CRP.Range("M55").ClearContents()
$ws_CRP.Range("N55").Value = -5685

# Row 107
$ws_CRP.Range("H107").Value = 568.5476
$ws_CRP.Range("I107").Value = 452.81482
$ws_CRP.Range("J107").Value = 776.86664
$ws_CRP.Range("K107").Value = 452.81482
$ws_CRP.Range("L107").Value = 776.86664
$ws_CRP.Range("M107").Value = 1467.18518
$ws_CRP.Range("N107").Value = -4616.86664

# Row 134
$ws_CRP.Range("H134").Value = 684.3823
$ws_CRP.Range("I134").Value = 679.65216
$ws_CRP.Range("J134").Value = 694.2727
$ws_CRP.Range("K134").Value = 2038.95648
$ws_CRP.Range("L134").Value = 2082.8181
$ws_CRP.Range("M134").Value = 496.0435200000002
$ws_CRP.Range("N134").Value = -7152.8181

$ws_CUL = $wb.Worksheets.Item("CUL")
# Row 5
$ws_CUL.Range("H5").Value = 1210.75
$ws_CUL.Range("I5").Value = 1155.0769
$ws_CUL.Range("J5").Value = 1276.5454
$ws_CUL.Range("K5").Value = 3465.2307
$ws_CUL.Range("L5").Value = 3829.6362
$ws_CUL.Range("M5").Value = -3353.2307
$ws_CUL.Range("N5").Value = -4053.6362

# Row 131
$ws_CUL.Range("H131").Value = 13514538
$ws_CUL.Range("I131").Value = 1690
$ws_CUL.Range("J131").Value = 16667536
$ws_CUL.Range("K131").Value = 5070
$ws_CUL.Range("L131").Value = 50002608
$ws_CUL.Range("M131").Value = -30
$ws_CUL.Range("N131").Value = -50012688

# Row 135
$ws_CUL.Range("H135").Value = 1210.75
$ws_CUL.Range("I135").Value = 1155.0769
$ws_CUL.Range("J135").Value = 1276.5454
$ws_CUL.Range("K135").Value = 10395.6921
$ws_CUL.Range("L135").Value = 11488.9086
$ws_CUL.Range("M135").Value = -7860.6921
$ws_CUL.Range("N135").Value = -16558.9086

$ws_GSM = $wb.Worksheets.Item("GSM")
# Row 52
$ws_GSM.Range("H52").Value = 0
$ws_GSM.Range("I52").Value = 0
$ws_GSM.Range("J52").Value = 0
$ws_GSM.Range("K52").Value = 0
$ws_GSM.Range("L52").ClearContents()
$ws_GSM.Range("N52").Value = 0

$ws_LTW = $wb.Worksheets.Item("LTW")
# Row 51
$ws_LTW.Range("H51").Value = 5163.364
$ws_LTW.Range("I51").Value = 0
$ws_LTW.Range("J51").Value = 5163.364
$ws_LTW.Range("K51").Value = 0
$ws_LTW.Range("L51").Value = 5163.364
$ws_LTW.Range("N51").Value = -6119.364

# Row 68
$ws_LTW.Range("H68").Value = 2126.6538
$ws_LTW.Range("I68").Value = 2016.6666
$ws_LTW.Range("J68").Value = 2374.125
$ws_LTW.Range("K68").Value = 2016.6666
$ws_LTW.Range("L68").Value = 2374.125
$ws_LTW.Range("M68").Value = -1267.6666
$ws_LTW.Range("N68").Value = -3872.125

# Row 71
$ws_LTW.Range("H71").Value = 2126.6538
$ws_LTW.Range("I71").Value = 2016.6666
$ws_LTW.Range("J71").Value = 2374.125
$ws_LTW.Range("K71").Value = 10083.333
$ws_LTW.Range("L71").Value = 11870.625
$ws_LTW.Range("M71").Value = -6339.333000000001
$ws_LTW.Range("N71").Value = -19358.625

$ws_WVR = $wb.Worksheets.Item("WVR")
# Row 62
$ws_WVR.Range("H62").Value = 5732.6665
$ws_WVR.Range("I62").Value = 5099
$ws_WVR.Range("J62").Value = 7000
$ws_WVR.Range("K62").Value = 5099
$ws_WVR.Range("L62").Value = 7000
$ws_WVR.Range("M62").Value = -4475
$ws_WVR.Range("N62").Value = -8248

# Row 65
$ws_WVR.Range("H65").Value = 5732.6665
$ws_WVR.Range("I65").Value = 5099
$ws_WVR.Range("J65").Value = 7000
$ws_WVR.Range("K65").Value = 25495
$ws_WVR.Range("L65").Value = 35000
$ws_WVR.Range("M65").Value = -22375
$ws_WVR.Range("N65").Value = -41240

# Row 81
$ws_WVR.Range("H81").Value = 66670340
$ws_WVR.Range("I81").Value = 3097.5454
$ws_WVR.Range("J81").Value = 250005250
$ws_WVR.Range("K81").Value = 6195.0908
$ws_WVR.Range("L81").Value = 500010500
$ws_WVR.Range("M81").Value = -5134.0908
$ws_WVR.Range("N81").Value = -500012622

# Row 84
$ws_WVR.Range("H84").Value = 66670340
$ws_WVR.Range("I84").Value = 3097.5454
$ws_WVR.Range("J84").Value = 250005250
$ws_WVR.Range("K84").Value = 30975.454
$ws_WVR.Range("L84").Value = 2500052500
$ws_WVR.Range("M84").Value = -25671.454
$ws_WVR.Range("N84").Value = -2500063108

# Row 92
$ws_WVR.Range("H92").Value = 22775
$ws_WVR.Range("I92").Value = 0
$ws_WVR.Range("J92").Value = 22775
$ws_WVR.Range("K92").Value = 0
$ws_WVR.Range("L92").Value = 22775
$ws_WVR.Range("N92").Value = -27767

# Row 136
$ws_WVR.Range("H136").Value = 976.2895
$ws_WVR.Range("I136").Value = 566.1724
$ws_WVR.Range("J136").Value = 2297.7778
$ws_WVR.Range("K136").Value = 1698.5172
$ws_WVR.Range("L136").Value = 2297.7778
$ws_WVR.Range("M136").Value = 851.4827999999998
$ws_WVR.Range("N136").Value = -11993.3334
